# Updates cryptos list values (price & volume columns, plus a couple of
# row re-orderings) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.932.50"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.642.09"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "'215.26"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'0.5083"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.2568"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.06394"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "'0.07772"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'4.298"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "1.653.05"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "'0.5452"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "0.0₅7857"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "'64.59"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "26.013.09"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'198.23"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "'4.435"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "'9.956"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "'6.036"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "'1.009"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'1.883"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "'139.77"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'0.1144"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'6.905"
$ws.Range("E27").Value = "  +3.04%  "
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'1.239"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'0.05005"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'3.262"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'3.186"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'1.537"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'2.362"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").Value = "'0.8921"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'2.592"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "1.129.81"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").Value = "'0.5509"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "'0.01560"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'1.004"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "'2.539"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'5.631"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₈126"
$ws.Range("E43").Value = "  +9.56%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8143"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'99.67"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "1.786.92"
$ws.Range("D47").Value = "'0.4525"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'54.82"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'0.05084"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  +0.29%  "
